$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new rows first (bottom-most insert first so row numbers
#     for the earlier insert are not disturbed) ---
# Old layout (rows 2-8):
#   2 U1        3 FUNC1      4 U2      5 S1       6 Rpot2     7 R1     8 R2
# New layout (rows 2-10):
#   2 U1 3 FUNC1 4 U2 5-6 S inicio/S captura 7 Rpot2 8 R1 9-10 R2/R3

# Insert a row after row 8 (R2) -> becomes row 9, for "R3"
$ws.Rows.Item(9).Insert()
# Insert a row after row 5 (S1) -> becomes row 6, for "S captura"
$ws.Rows.Item(6).Insert()

# --- Copy down formatting for the two freshly inserted blank rows ---
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:C9").Copy()
$ws.Range("A10:C10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 3 (FUNC1): update the function-generator description ---
$ws.Range("C3").Value = "1 Hz, 2 V, 0 V, Seno Generador de función"

# --- Rows 5/6 (was "S1" row): split into "S inicio" / "S captura" ---
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "S inicio"
$ws.Range("C5").Value = " Pulsador"
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = "S captura"
$ws.Range("C6").Value = ""

# --- Rows 9/10 (was "R2" row): split into "R2" / "R3" ---
$ws.Range("A9").Value = 2
$ws.Range("B9").Value = "R2"
$ws.Range("C9").Value = "10 kΩ Resistencia"
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = "R3"
$ws.Range("C10").Value = ""

# --- Merge the split rows ---
$ws.Range("A5:A6").Merge() | Out-Null
$ws.Range("C5:C6").Merge() | Out-Null
$ws.Range("A9:A10").Merge() | Out-Null
$ws.Range("C9:C10").Merge() | Out-Null

# --- Wrap text on the "S captura" cell ---
$ws.Range("B6").WrapText = $true

# --- Drop the custom 15.75pt row height on the data rows (back to default) ---
$ws.Range("A2:C10").Rows.AutoFit() | Out-Null

# --- Selection / view state ---
$ws.Range("D2").Select() | Out-Null

Write-Host "done"
